# Sampling_Resolution_Calculator.xlsx — "added beam diagnostics for parameter
# scans. f1_scan submitted at NERSC"
#
# All the changed cells live on the "HRM" worksheet. The large cascade of
# <v> (cached value) changes in the diff all fall out automatically from a
# handful of root-cause edits once the workbook recalculates, so we only
# need to touch the inputs/formulas below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HRM")

# --- Root-cause input / formula edits -------------------------------------

# B3: 100 -> 400
$ws.Range("B3").Value = 400

# B4 / B5: literal 4.0E-3 -> formula =0.005 (evaluates to 5.0E-3)
$ws.Range("B4").Formula = "=0.005"
$ws.Range("B5").Formula = "=0.005"

# I4: MAX(ROUND(I3/$B7,0),1) -> MAX(ROUND(2*I3/$B7,0),1)
$ws.Range("I4").Formula = "=MAX(ROUND(2*I3/`$B7,0),1)"

# B10: MIN(B9/400,0.001) -> MIN(B9/800,0.001)
$ws.Range("B10").Formula = "=MIN(B9/800,0.001)"

# B12: 512 -> 1024
$ws.Range("B12").Value = 1024

# F22: D22*I2*4*4 -> D22*I2*8
$ws.Range("F22").Formula = "=D22*I2*8"

# P24: N24*2 -> N24
$ws.Range("P24").Formula = "=N24"

# --- Column width tweaks ----------------------------------------------------
# col I (9): 8.375 -> 9.375 ; col N (14): 9.375 -> 10.375 ; col P (16): 9.375 -> 10.375
# (ColumnWidth is in "characters"; the stored sheet XML "width" attribute is
# ColumnWidth + 5/7, so we back that offset out of the target stored width.)
$ws.Columns.Item(9).ColumnWidth = 9.375 - 5/7
$ws.Columns.Item(14).ColumnWidth = 10.375 - 5/7
$ws.Columns.Item(16).ColumnWidth = 10.375 - 5/7

# --- Selection -------------------------------------------------------------
$ws.Activate()
$ws.Range("L36").Select()

$wb.Application.Calculate()
